# Apply the update described by the diff:
#  - Row 15: clear the (already empty) placeholder cells F15, G15 and L15 so
#    they no longer exist as explicit cells.
#  - Append a brand-new row 16 with the data for the new order, including
#    empty placeholder cells for Optimizador (F), Unidades Optimizador (G),
#    Baterías (J), Unidades Baterías (K) and Cargador VE (L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty placeholder cells left over on row 15.
$ws.Range("F15").ClearContents()
$ws.Range("G15").ClearContents()
$ws.Range("L15").ClearContents()

# New row 16 data. Numero (A) is a real number; Unidades columns (E, I) are
# stored as *text* digits in the source data (same convention as the rest of
# the sheet), so they are written with a leading apostrophe to keep them as
# text instead of being auto-coerced to numbers.
$ws.Range("A16").Value = 1598
$ws.Range("B16").Value = "Juan José Lopez"
$ws.Range("C16").Value = "Estructura coplanar NOVOTEGRA"
$ws.Range("D16").Value = "Trina 505W TSM-NEG18R.25"
$ws.Range("E16").Value = "'19"
$ws.Range("H16").Value = "GOODWE ES UNIQ - GW10000-ES-C10 híbrido monofásico"
$ws.Range("I16").Value = "'1"
$ws.Range("M16").Value = "Sí"
$ws.Range("N16").Value = "2024-01-03T10:49:29.088Z"
